# Updates the crypto price/volume table to reflect the latest scrape.
# Column D (Price) values are assigned with a leading apostrophe so Excel
# stores them as text (preserving exact formatting such as trailing
# zeros or "thousands.dot" notation) instead of silently converting
# numeric-looking strings into numbers, matching the source data which
# is plain text in every cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'48.152.92"
$ws.Range("E2").Value = "  +1.73%  "

$ws.Range("D3").Value = "'2.509.73"
$ws.Range("E3").Value = "  +0.80%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'320.87"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").Value = "'108.83"
$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("E7").Value = "  +1.28%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  +1.00%  "

$ws.Range("D10").Value = "'39.96"
$ws.Range("E10").Value = "  +1.35%  "

$ws.Range("D11").Value = "'20.31"
$ws.Range("E11").Value = "  +10.68%  "

$ws.Range("D12").Value = "'0.0818"
$ws.Range("E12").Value = "  +0.88%  "

$ws.Range("E13").Value = "  +0.66%  "

$ws.Range("D14").Value = "'7.20"
$ws.Range("E14").Value = "  +0.83%  "

$ws.Range("D15").Value = "'2.902.31"
$ws.Range("E15").Value = "  +0.82%  "

$ws.Range("D16").Value = "'2.521.05"
$ws.Range("E16").Value = "  +0.95%  "

$ws.Range("D17").Value = "'0.848"
$ws.Range("E17").Value = "  +0.34%  "

$ws.Range("D18").Value = "'48.011.45"
$ws.Range("E18").Value = "  +1.64%  "

$ws.Range("D19").Value = "'13.18"
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("E21").Value = "  +0.75%  "

$ws.Range("E22").Value = "  +2.58%  "

$ws.Range("D23").Value = "'72.25"

$ws.Range("D24").Value = "'276.68"
$ws.Range("E24").Value = "  +12.78%  "

$ws.Range("D25").Value = "'2.57"
$ws.Range("E25").Value = "  +0.37%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").Value = "'25.91"
$ws.Range("E27").Value = "  +0.75%  "

$ws.Range("D28").Value = "'2.40"
$ws.Range("E28").Value = "  +5.59%  "

$ws.Range("D29").Value = "'10.06"
$ws.Range("E29").Value = "  +0.66%  "

$ws.Range("D30").Value = "'0.141"
$ws.Range("E30").Value = "  +2.16%  "

$ws.Range("D31").Value = "'35.40"
$ws.Range("E31").Value = "  +2.06%  "

$ws.Range("D32").Value = "'49.54"

$ws.Range("D33").Value = "'19.32"
$ws.Range("E33").Value = "  -6.05%  "

$ws.Range("D34").Value = "'5.35"
$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("D36").Value = "'0.0786"
$ws.Range("E36").Value = "  +0.19%  "

$ws.Range("E37").Value = "  -0.43%  "

$ws.Range("E38").Value = "  -2.11%  "

$ws.Range("E39").Value = "  +1.16%  "

# Rows 40 and 41 swap places: Monero (previously row 40) moves to row 41
# and Stellar (previously row 41) moves to row 40, each with refreshed
# price / volume figures.
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "'0.112"
$ws.Range("E40").Value = "  +0.49%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'122.58"
$ws.Range("E41").Value = "  +5.11%  "

$ws.Range("D42").Value = "'2.22"
$ws.Range("E42").Value = "  -0.75%  "

$ws.Range("D43").Value = "'21.70"
$ws.Range("E43").Value = "  -6.89%  "

$ws.Range("D44").Value = "'0.0307"
$ws.Range("E44").Value = "  +3.53%  "

$ws.Range("D45").Value = "'2.019.87"
$ws.Range("E45").Value = "  +1.10%  "

$ws.Range("D46").Value = "'3.13"
$ws.Range("E46").Value = "  +3.03%  "

$ws.Range("E47").Value = "  +4.46%  "

$ws.Range("E48").Value = "  -0.35%  "

$ws.Range("D49").Value = "'9.05"
$ws.Range("E49").Value = "  -1.34%  "

$ws.Range("E50").Value = "  +1.99%  "

$ws.Range("D51").Value = "'79.75"
$ws.Range("E51").Value = "  +2.86%  "
